# Weekly fruit/vegetable price update: a new row of price data is inserted
# at row 485 (pushing the existing rows 485-523 down to 486-524), matching
# the canonical diff for this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 485, shifting rows 485:523 down to 486:524.
$ws.Rows.Item(485).Insert()

# Populate the newly inserted row 485 with the new weekly record.
$ws.Cells.Item(485, 1).Value = 9
$ws.Cells.Item(485, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(485, 3).Value = "Metropolitana"
$ws.Cells.Item(485, 4).Value = 45265
$ws.Cells.Item(485, 5).Value = 13
$ws.Cells.Item(485, 6).Value = 300000001
$ws.Cells.Item(485, 7).Value = "Rabanito"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 7000
$ws.Cells.Item(485, 11).Value = 3000
$ws.Cells.Item(485, 12).Value = 3000
$ws.Cells.Item(485, 13).Value = 3000
$ws.Cells.Item(485, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(485, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(485, 16).Value = 30
$ws.Cells.Item(485, 17).Value = 100
$ws.Cells.Item(485, 18).Value = "Hortaliza"
